$d = $word.ActiveDocument

# Change 1: Simplify the hardware design / UART / I2C sentence
$d.Content.Find.Execute(
    "I have knowledge about the hardware design and the UART, I2C protocols for sending/receiving data ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I have knowledge about protocols for sending/receiving data ",
    2) | Out-Null

# Change 2: Add "and automation test" after "IT test"
$d.Content.Find.Execute(
    ", writing UT test, IT test for ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", writing UT test, IT test and automation test for ",
    2) | Out-Null

# Change 3: Add "and Qt Framework as well" after "CAN bus protocol"
$d.Content.Find.Execute(
    "CAN bus protocol. I use",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CAN bus protocol and Qt Framework as well. I use",
    2) | Out-Null

# Change 4: Add ", Linux command line" after "git, SVN, Visual Code"
$d.Content.Find.Execute(
    "git, SVN, Visual Code ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "git, SVN, Visual Code, Linux command line ",
    2) | Out-Null
